$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44495
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 11000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 11500
$ws.Range("S2").Value = 5750

$ws.Range("D3").Value = 44459
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 13000
$ws.Range("O3").Value = 14000
$ws.Range("P3").Value = 13500
$ws.Range("S3").Value = 6750

$ws.Range("D4").Value = 44491
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 11500
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 11750
$ws.Range("S4").Value = 5875

$ws.Range("D5").Value = 44490
$ws.Range("M5").Value = 160
$ws.Range("N5").Value = 11500
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 11750
$ws.Range("S5").Value = 5875

$ws.Range("D6").Value = 44448
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14500
$ws.Range("S6").Value = 7250

$ws.Range("D7").Value = 44497
$ws.Range("M7").Value = 400
$ws.Range("N7").Value = 11500
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 11750
$ws.Range("S7").Value = 5875

$ws.Range("D8").Value = 44446
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("S8").Value = 7250

$ws.Range("D9").Value = 44498
$ws.Range("M9").Value = 240
$ws.Range("N9").Value = 11000
$ws.Range("O9").Value = 11500
$ws.Range("P9").Value = 11250
$ws.Range("S9").Value = 5625

$ws.Range("D10").Value = 44455
$ws.Range("M10").Value = 160
$ws.Range("N10").Value = 13000
$ws.Range("O10").Value = 14000
$ws.Range("P10").Value = 13500
$ws.Range("S10").Value = 6750

$ws.Range("D11").Value = 44466
$ws.Range("M11").Value = 160
$ws.Range("N11").Value = 13500
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 13750
$ws.Range("S11").Value = 6875

$ws.Range("D12").Value = 44462
$ws.Range("M12").Value = 140
$ws.Range("N12").Value = 13000
$ws.Range("O12").Value = 14000
$ws.Range("P12").Value = 13500
$ws.Range("S12").Value = 6750

$ws.Range("D13").Value = 44463
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 13000
$ws.Range("O13").Value = 14000
$ws.Range("P13").Value = 13500
$ws.Range("S13").Value = 6750

$ws.Range("D14").Value = 44452
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 13000
$ws.Range("O14").Value = 14000
$ws.Range("P14").Value = 13500
$ws.Range("S14").Value = 6750

$ws.Range("D15").Value = 44454
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 13000
$ws.Range("O15").Value = 14000
$ws.Range("P15").Value = 13500
$ws.Range("S15").Value = 6750

$ws.Range("D16").Value = 44445
$ws.Range("M16").Value = 160
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 14500
$ws.Range("S16").Value = 7250

$ws.Range("D17").Value = 44494
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 11500
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 11750
$ws.Range("S17").Value = 5875

$ws.Range("D18").Value = 44489
$ws.Range("M18").Value = 400
$ws.Range("N18").Value = 11500
$ws.Range("O18").Value = 12000
$ws.Range("P18").Value = 11750
$ws.Range("S18").Value = 5875

$ws.Range("D19").Value = 44468
$ws.Range("M19").Value = 300
$ws.Range("N19").Value = 13000
$ws.Range("O19").Value = 14000
$ws.Range("P19").Value = 13500
$ws.Range("S19").Value = 6750
